# Revert capacity chart to show kilowatts on the y-axis.
# The underlying data was previously recorded in Watts; it now needs to
# be expressed in kilowatts (divide by 1000), the number format gets an
# extra decimal place, and the axis title/number format switch from the
# old "Watts" / "K-suffix" presentation back to plain kilowatts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Worksheet values: convert Watts -> kilowatts -----------------
$ws.Range("C25").Value = 3.84

$ws.Range("E12").Value = 75.30000000000001
$ws.Range("E14").Value = 5.58
$ws.Range("E15").Value = 5.8
$ws.Range("E16").Value = 13
$ws.Range("E17").Value = 12.3
$ws.Range("E18").Value = 33.45
$ws.Range("E19").Value = 67.5
$ws.Range("E20").Value = 31.9
$ws.Range("E21").Value = 42.1
$ws.Range("E22").Value = 52.8
$ws.Range("E23").Value = 106.024
$ws.Range("E24").Value = 310.33
$ws.Range("E25").Value = 333.34
$ws.Range("E26").Value = 53.532

# --- 2. Number format: show one decimal place now that values are kW -
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- 3. Chart: axis title + value-axis number format ------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
